$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.329283714294434
$ws.Range("B1").Value = 2.473553657531738
$ws.Range("C1").Value = 2.056790590286255
$ws.Range("D1").Value = 1.948425531387329
$ws.Range("E1").Value = 1.714047193527222
